$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.358.32'
$ws.Range('E2').Value = '  +0.69%  '
$ws.Range('D3').Value = '1.849.16'
$ws.Range('E3').Value = '  +0.67%  '
$ws.Range('D4').Value = "'1.014"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +1.05%  '
$ws.Range('D5').Value = "'244.62"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.11%  '
$ws.Range('D6').Value = "'0.6197"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.63%  '
$ws.Range('D7').Value = "'1.013"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.72%  '
$ws.Range('D8').Value = "'0.07466"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.16%  '
$ws.Range('D9').Value = "'0.2958"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.86%  '
$ws.Range('E10').Value = '  -0.12%  '
$ws.Range('D11').Value = "'0.07742"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.00%  '
$ws.Range('D12').Value = '1.843.40'
$ws.Range('E12').Value = '  +0.40%  '
$ws.Range('D13').Value = "'5.022"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.25%  '
$ws.Range('D14').Value = "'0.6746"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.76%  '
$ws.Range('D15').Value = "'83.37"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.24%  '
$ws.Range('D16').Value = "'0.000009081"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -3.49%  '
$ws.Range('D17').Value = "'5.908"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.63%  '
$ws.Range('D18').Value = '29.327.15'
$ws.Range('E18').Value = '  +0.50%  '
$ws.Range('D19').Value = '2.085.58'
$ws.Range('E19').Value = '  +0.17%  '
$ws.Range('D20').Value = "'238.33"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +6.40%  '
$ws.Range('E21').Value = '  +0.60%  '
$ws.Range('E22').Value = '  +0.77%  '
$ws.Range('D23').Value = "'7.197"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.89%  '
$ws.Range('D24').Value = "'1.016"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.18%  '
$ws.Range('D25').Value = "'160.11"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.22%  '
$ws.Range('D26').Value = "'0.1435"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.14%  '
$ws.Range('D27').Value = "'8.544"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.32%  '
$ws.Range('D28').Value = "'17.95"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.06%  '
$ws.Range('D29').Value = "'1.507"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.18%  '
$ws.Range('D30').Value = "'4.170"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.43%  '
$ws.Range('D31').Value = "'0.05610"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.39%  '
$ws.Range('D32').Value = "'4.123"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.31%  '
$ws.Range('D33').Value = "'1.222"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.61%  '
$ws.Range('D34').Value = "'0.7530"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.01%  '
$ws.Range('D35').Value = "'1.857"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.05%  '
$ws.Range('D36').Value = "'1.146"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.61%  '
$ws.Range('D37').Value = "'2.675"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.16%  '
$ws.Range('D38').Value = "'2.829"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.52%  '
$ws.Range('D39').Value = "'0.01789"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.14%  '
$ws.Range('D40').Value = '1.216.72'
$ws.Range('E40').Value = '  -1.09%  '
$ws.Range('D41').Value = "'6.511"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.35%  '
$ws.Range('D42').Value = "'0.9046"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.92%  '
$ws.Range('D43').Value = "'1.013"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.62%  '
$ws.Range('B44').Value = 'Quant'
$ws.Range('C44').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D44').Value = "'101.45"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.47%  '
$ws.Range('B45').Value = 'RocketPoolETH'
$ws.Range('C45').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D45').Value = '1.992.77'
$ws.Range('E45').Value = '  +0.78%  '
$ws.Range('D46').Value = "'65.56"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.12%  '
$ws.Range('E47').Value = '  -1.13%  '
$ws.Range('D48').Value = "'0.5152"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.98%  '
$ws.Range('B49').Value = 'TheSandbox'
$ws.Range('C49').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D49').Value = "'0.4072"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.38%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = "'9.190"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.36%  '
$ws.Range('D51').Value = "'0.05852"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.80%  '
